$wb = $excel.ActiveWorkbook

# --- config_netConnections: re-point/aggregate owner_actor (col I) for rows 21-25 ---
# Before: com1, com1, com1, com2, com2
# After : com2, com3, com3, com4, com4  (com3/com4 are brand-new aggregate actor types)
$wsConn = $wb.Worksheets.Item("config_netConnections")
$wsConn.Range("I21").Value = "com2"
$wsConn.Range("I22").Value = "com3"
$wsConn.Range("I23").Value = "com3"
$wsConn.Range("I24").Value = "com4"
$wsConn.Range("I25").Value = "com4"

# --- Update the saved cursor/selection + active worksheet tab ---
# config_netConnections is no longer the selected tab; leave its cursor at N27
$wsConn.Activate() | Out-Null
$wsConn.Range("N27").Select() | Out-Null

# config_actors becomes the selected/active tab; leave its cursor at E20
$wsActors = $wb.Worksheets.Item("config_actors")
$wsActors.Activate() | Out-Null
$wsActors.Range("E20").Select() | Out-Null
